# Daily update at 8 AM UTC: append the next day's row of data and
# move the "last row" date-only formatting down to the new last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The previous last row (35) had a terminal "date only" number format;
# now that it's no longer the last row it goes back to the standard
# "date + time" format used by every other data row.
$ws.Range("A35").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new day's data as row 36.
$ws.Range("A36").Value = 45776
$ws.Range("B36").Value = 145
$ws.Range("C36").Value = 150
$ws.Range("D36").Value = 147

# The new last row gets the terminal "date only" number format.
$ws.Range("A36").NumberFormat = "YYYY-MM-DD"
